$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ERSL-Test Cases")

# Row 6 (ERSL-4: Check login with invalid username) - new test data
$ws.Range("D6").Value = "Username = bipul51`nPassword  = password"

# Row 7 (ERSL-6: Check login with invalid password) - new expected/actual results first
$ws.Range("E7").Value = "Display error message as:`n`"Username and Password does not match. Please try again with correct credentials !!!`" "
$ws.Range("F7").Value = "Username and Password does not match.`r`n Please try again with correct credentials !!!"

# then new test data for row 7
$ws.Range("D7").Value = "Username = abc1234`nPassword = testing123"

# Remaining row 6 cells (reuse existing shared strings)
$ws.Range("E6").Value = "Display error message as:`n`"Username and Password does not match. Please try again with correct credentials !!!`""
$ws.Range("F6").Value = "`"Username and Password does not match.`n Please try again with correct credentials !!!`""
$ws.Range("G6").Value = "Pass"
$ws.Range("G7").Value = "Pass"

$ws.Range("E6").WrapText = $true
$ws.Range("F6").WrapText = $true
$ws.Range("G6").WrapText = $true
$ws.Range("E7").WrapText = $true
$ws.Range("F7").WrapText = $true
$ws.Range("G7").WrapText = $true

# Match the author's final UI selection
$ws.Range("D7").Select()
